$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3205216666666666
$ws.Range("H2").Value = 0.961565
$ws.Range("I2").Value = 0.1355444593128396
$ws.Range("J2").Value = 0.1355444593128396
$ws.Range("M2").Value = 2.031394333333334
$ws.Range("N2").Value = 6.094183
$ws.Range("O2").Value = 0.1230522080665953
$ws.Range("P2").Value = 0.1230522080665953
$ws.Range("Q2").Value = 0.6511058973772222
$ws.Range("R2").Value = 5.859953076395001
$ws.Range("S2").Value = 0.01667904500963769
$ws.Range("T2").Value = 0.01667904500963769

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3205216666666666
$ws.Range("H3").Value = 0.961565
$ws.Range("I3").Value = 0.1355444593128396
$ws.Range("J3").Value = 0.1355444593128396
$ws.Range("O3").Value = 0.6442150581092261
$ws.Range("P3").Value = 0.6442150581092261
$ws.Range("Q3").Value = 3.408733822046667
$ws.Range("R3").Value = 30.67860439842
$ws.Range("S3").Value = 0.08731978173260457
$ws.Range("T3").Value = 0.08731978173260457

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3205216666666666
$ws.Range("H4").Value = 0.961565
$ws.Range("I4").Value = 0.1355444593128396
$ws.Range("J4").Value = 0.1355444593128396
$ws.Range("M4").Value = 3.842043666666667
$ws.Range("N4").Value = 11.526131
$ws.Range("O4").Value = 0.2327327338241786
$ws.Range("P4").Value = 0.2327327338241786
$ws.Range("Q4").Value = 1.231458239446111
$ws.Range("R4").Value = 11.083124155015
$ws.Range("S4").Value = 0.03154563257059729
$ws.Range("T4").Value = 0.0315456325705973

# Row 5
$ws.Range("I5").Value = 0.4633580360449179
$ws.Range("J5").Value = 0.4633580360449179
$ws.Range("M5").Value = 2.031394333333334
$ws.Range("N5").Value = 6.094183
$ws.Range("O5").Value = 0.1230522080665953
$ws.Range("P5").Value = 0.1230522080665953
$ws.Range("Q5").Value = 2.225802156690555
$ws.Range("R5").Value = 20.032219410215
$ws.Range("S5").Value = 0.05701722946072819
$ws.Range("T5").Value = 0.0570172294607282

# Row 6
$ws.Range("I6").Value = 0.4633580360449179
$ws.Range("J6").Value = 0.4633580360449179
$ws.Range("O6").Value = 0.6442150581092261
$ws.Range("P6").Value = 0.6442150581092261
$ws.Range("S6").Value = 0.2985022241160537
$ws.Range("T6").Value = 0.2985022241160537

# Row 7
$ws.Range("I7").Value = 0.4633580360449179
$ws.Range("J7").Value = 0.4633580360449179
$ws.Range("M7").Value = 3.842043666666667
$ws.Range("N7").Value = 11.526131
$ws.Range("O7").Value = 0.2327327338241786
$ws.Range("P7").Value = 0.2327327338241786
$ws.Range("Q7").Value = 4.209733648972778
$ws.Range("R7").Value = 37.887602840755
$ws.Range("S7").Value = 0.107838582468136
$ws.Range("T7").Value = 0.107838582468136

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9484743333333334
$ws.Range("H8").Value = 2.845423
$ws.Range("I8").Value = 0.4010975046422426
$ws.Range("J8").Value = 0.4010975046422425
$ws.Range("M8").Value = 2.031394333333334
$ws.Range("N8").Value = 6.094183
$ws.Range("O8").Value = 0.1230522080665953
$ws.Range("P8").Value = 0.1230522080665953
$ws.Range("Q8").Value = 1.926725386045445
$ws.Range("R8").Value = 17.340528474409
$ws.Range("S8").Value = 0.0493559335962294
$ws.Range("T8").Value = 0.0493559335962294

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9484743333333334
$ws.Range("H9").Value = 2.845423
$ws.Range("I9").Value = 0.4010975046422426
$ws.Range("J9").Value = 0.4010975046422425
$ws.Range("O9").Value = 0.6442150581092261
$ws.Range("P9").Value = 0.6442150581092261
$ws.Range("Q9").Value = 10.08698280212933
$ws.Range("R9").Value = 90.782845219164
$ws.Range("S9").Value = 0.2583930522605679
$ws.Range("T9").Value = 0.2583930522605679

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9484743333333334
$ws.Range("H10").Value = 2.845423
$ws.Range("I10").Value = 0.4010975046422426
$ws.Range("J10").Value = 0.4010975046422425
$ws.Range("M10").Value = 3.842043666666667
$ws.Range("N10").Value = 11.526131
$ws.Range("O10").Value = 0.2327327338241786
$ws.Range("P10").Value = 0.2327327338241786
$ws.Range("Q10").Value = 3.644079805379223
$ws.Range("R10").Value = 32.79671824841301
$ws.Range("S10").Value = 0.09334851878544527
$ws.Range("T10").Value = 0.09334851878544527

